# Add 2 APIs: "Get Call History by user" and "Get Unique Dispositions"
# (plus a third "Dashboard Voice Bot" entry that is also present in the target
# sheet) as three new rows appended after the existing API table, each
# separated by a blank row like the rest of the sheet (rows 57, 60, 63).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 57: Dashboard Voice Bot ( GET ) ----------------------------------
$ws.Range("A57").Value = "Dashboard Voice Bot ( GET )"
$b57 = $ws.Range("B57")
$b57.Hyperlinks.Add($b57, "http://1msg.1point1.in:3001/api/auth/j-v1/dashboard/voice/bot/?user_id=6") | Out-Null
$ws.Range("C57").Value = "user_id=6"
$ws.Range("D57").Value = "curl --location 'http://1msg.1point1.in:3001/api/auth/j-v1/dashboard/voice/bot/?user_id=6'"

# --- Row 60: Get Call History by Agent ID ( GET ) -------------------------
$ws.Range("A60").Value = "Get Call History by Agent ID    ( GET )   "
$b60 = $ws.Range("B60")
$b60.Hyperlinks.Add($b60, "http://1msg.1point1.in:3001/api/auth/j-v1/call-history/by/user_id/?user_id=1") | Out-Null
$ws.Range("C60").Value = "user_id=10"
$ws.Range("D60").Value = "curl --location 'http://1msg.1point1.in:3001/api/auth/j-v1/call-history/by/user_id/?user_id=10'"

# --- Row 63: Get Unique Dispositions ( GET ) -------------------------------
$ws.Range("A63").Value = "Get Unique Dispositions ( GET )  "
$b63 = $ws.Range("B63")
$b63.Hyperlinks.Add($b63, "http://1msg.1point1.in:3001/api/auth/j-v1/get/unique-disposition/") | Out-Null
$ws.Range("D63").Value = "curl --location 'http://1msg.1point1.in:3001/api/auth/j-v1/get/unique-disposition/'"

# Hyperlinks.Add() always mints a brand-new cell style; put the three new
# link cells back on the workbook's existing shared "Hyperlink" style so
# they match the other hyperlinked cells in the sheet (e.g. B45, B54).
$b57.Style = "Hyperlink"
$b60.Style = "Hyperlink"
$b63.Style = "Hyperlink"

# Leave the selection on the last new cell, matching where the author was
# working when the rows were added.
$ws.Range("D63").Select()
